$d = $word.ActiveDocument

# The document ends with a trailing empty "Normal (Web)" paragraph right
# before the section break. We need to add 6 new paragraphs *before* that
# existing trailing paragraph (it must stay last, untouched):
#   - 4 blank "Normal (Web)" paragraphs (same formatting as the trailing one)
#   - a paragraph containing "{{ izjava }}"
#   - a tab-indented paragraph containing "{{ potpis }}"
#
# Plain Range.InsertParagraphAfter()/Before() calls leave behind a stray
# empty <w:r> run carrying the current character formatting on every newly
# minted paragraph, which the target markup does not have. To avoid that,
# we first grow the document by one new (still-empty) paragraph placed
# immediately before the existing trailing paragraph, then overwrite that
# new paragraph's own range with the exact WordprocessingML for all six
# paragraphs via Range.InsertXML - this replaces the placeholder run
# instead of appending to it, so the inserted markup matches exactly.

$lastIndex = $d.Paragraphs.Count
$trailingPara = $d.Paragraphs($lastIndex)

$trailingPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($lastIndex)
$newRange = $newPara.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0" w:line="216" w:lineRule="atLeast"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0" w:line="216" w:lineRule="atLeast"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0" w:line="216" w:lineRule="atLeast"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0" w:line="216" w:lineRule="atLeast"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0" w:line="216" w:lineRule="atLeast"/><w:ind w:firstLine="284"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>{{ izjava }}</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="7088"/></w:tabs><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">{{ </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>potpis</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve"> }}</w:t></w:r></w:p>
'@

$newRange.InsertXML($xml)
